$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.854.31"
$ws.Range("E2").Value = "  +0.43%  "

Set-TextValue $ws.Range("D3") "1.857.68"
$ws.Range("E3").Value = "  -0.37%  "

Set-TextValue $ws.Range("D4") "1.016"
$ws.Range("E4").Value = "  -1.94%  "

Set-TextValue $ws.Range("D5") "321.04"
$ws.Range("E5").Value = "  -1.11%  "

Set-TextValue $ws.Range("D6") "1.015"
$ws.Range("E6").Value = "  -1.85%  "

Set-TextValue $ws.Range("D7") "0.4323"
$ws.Range("E7").Value = "  -2.06%  "

Set-TextValue $ws.Range("D8") "0.3801"

Set-TextValue $ws.Range("D9") "0.07420"
$ws.Range("E9").Value = "  -0.56%  "

Set-TextValue $ws.Range("D10") "0.8857"
$ws.Range("E10").Value = "  +0.14%  "

Set-TextValue $ws.Range("D11") "21.78"
$ws.Range("E11").Value = "  +0.03%  "

Set-TextValue $ws.Range("D12") "1.871.22"
$ws.Range("E12").Value = "  +0.12%  "

Set-TextValue $ws.Range("D13") "6.780"
$ws.Range("E13").Value = "  +0.39%  "

Set-TextValue $ws.Range("D14") "5.503"
$ws.Range("E14").Value = "  -0.98%  "

Set-TextValue $ws.Range("D15") "0.07108"
$ws.Range("E15").Value = "  -1.55%  "

Set-TextValue $ws.Range("D16") "88.77"
$ws.Range("E16").Value = "  +5.79%  "

$ws.Range("E17").Value = "  -1.84%  "

Set-TextValue $ws.Range("D18") "0.000009061"
$ws.Range("E18").Value = "  -0.44%  "

Set-TextValue $ws.Range("D19") "1.015"
$ws.Range("E19").Value = "  -1.81%  "

Set-TextValue $ws.Range("D20") "15.56"
$ws.Range("E20").Value = "  +0.35%  "

Set-TextValue $ws.Range("D21") "27.890.31"
$ws.Range("E21").Value = "  +0.46%  "

Set-TextValue $ws.Range("D22") "5.288"
$ws.Range("E22").Value = "  -0.42%  "

Set-TextValue $ws.Range("D23") "11.21"
$ws.Range("E23").Value = "  -1.79%  "

Set-TextValue $ws.Range("D24") "2.095.69"
$ws.Range("E24").Value = "  +0.47%  "

Set-TextValue $ws.Range("D25") "2.032"
$ws.Range("E25").Value = "  +3.68%  "

Set-TextValue $ws.Range("D26") "156.81"
$ws.Range("E26").Value = "  -0.91%  "

Set-TextValue $ws.Range("D27") "18.70"
$ws.Range("E27").Value = "  -0.94%  "

Set-TextValue $ws.Range("D28") "2.028"
$ws.Range("E28").Value = "  +1.33%  "

Set-TextValue $ws.Range("D29") "5.428"
$ws.Range("E29").Value = "  +2.27%  "

Set-TextValue $ws.Range("D30") "121.90"
$ws.Range("E30").Value = "  +3.66%  "

Set-TextValue $ws.Range("D31") "0.08982"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("E32").Value = "  +2.72%  "

Set-TextValue $ws.Range("D33") "0.7778"
$ws.Range("E33").Value = "  +0.92%  "

Set-TextValue $ws.Range("D34") "4.588"
$ws.Range("E34").Value = "  +0.11%  "

Set-TextValue $ws.Range("D35") "2.936"
$ws.Range("E35").Value = "  -3.42%  "

Set-TextValue $ws.Range("D36") "1.150"
$ws.Range("E36").Value = "  -1.42%  "

Set-TextValue $ws.Range("D37") "1.016"
$ws.Range("E37").Value = "  -1.94%  "

Set-TextValue $ws.Range("D38") "0.05332"
$ws.Range("E38").Value = "  -0.43%  "

Set-TextValue $ws.Range("D39") "0.01972"
$ws.Range("E39").Value = "  -0.85%  "

Set-TextValue $ws.Range("D40") "2.886"
$ws.Range("E40").Value = "  +1.42%  "

Set-TextValue $ws.Range("D41") "0.5210"
$ws.Range("E41").Value = "  +0.25%  "

Set-TextValue $ws.Range("D42") "7.029"
$ws.Range("E42").Value = "  +2.58%  "

Set-TextValue $ws.Range("D43") "0.1685"
$ws.Range("E43").Value = "  -0.43%  "

Set-TextValue $ws.Range("D44") "8.811"
$ws.Range("E44").Value = "  +1.12%  "

Set-TextValue $ws.Range("D45") "110.88"
$ws.Range("E45").Value = "  +1.15%  "

Set-TextValue $ws.Range("D46") "10.71"
$ws.Range("E46").Value = "  +0.74%  "

Set-TextValue $ws.Range("D49") "0.06529"
$ws.Range("E49").Value = "  +1.54%  "

Set-TextValue $ws.Range("D50") "1.016"
$ws.Range("E50").Value = "  -2.00%  "

Set-TextValue $ws.Range("D51") "1.882"
$ws.Range("E51").Value = "  +0.50%  "

# Row 47 and 48 swap (NEARProtocol <-> Decentraland) with updated values
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.4762"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.717"
$ws.Range("E48").Value = "  -0.60%  "
